# "Add files via upload" — adds a new file-contribution row (row 36) for
# ClassRelationshipDiagram.draw.io to the "File Contributions" table, with
# 0% Design, 100% Documentation, 0% Implementation, 0% Testing, and a Total
# (SUM) formula consistent with the rest of the F column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New file-contribution entry.
$ws.Range("A36").Value = "ClassRelationshipDiagram.draw.io"
$ws.Range("B36").Value = 0
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("F36").Formula = "=SUM(B36:E36)"

# Match the author's final selection/scroll position from the saved file.
$ws.Range("F36").Select()
